$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# New rows of data to append (rows 12-21), matching the shared-string values
# used by the original commit (column order: Time, RunningTime(s), Preprocess,
# Features, Model, Model_Details, Test_Accuracy, Val_Accuracy, Template Filter,
# <blank col I is Template Filter>, <col J>)

$preprocessA = 'remove multiple spaces, trim "space" and ",", convert to lower, convert unicode to ascii'
$modelDetailsA = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000'

$preprocessB = 'trim "space" and ",", convert to lower, remove multiple spaces, convert unicode to ascii'
$modelDetailsB = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000'

$features = '7 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type'
$model = 'Neuron Network'
$templateFilter = '0 filters: '

$rows = @(
    @{ Row = 12; Time = '20160406_104448'; RunningTime = 1134.701; Preprocess = $preprocessA; ModelDetails = $modelDetailsA; Test = 0.999333333333333; Val = 0.867986798679868; J = 0.0333333333333333 },
    @{ Row = 13; Time = '20160406_110343'; RunningTime = 1138.221; Preprocess = $preprocessA; ModelDetails = $modelDetailsA; Test = 0.998666666666667; Val = 0.864686468646865; J = 0.0862068965517241 },
    @{ Row = 14; Time = '20160406_112241'; RunningTime = 1201.19;  Preprocess = $preprocessA; ModelDetails = $modelDetailsA; Test = 1;                 Val = 0.897689768976898; J = 0.101449275362319 },
    @{ Row = 15; Time = '20160406_114242'; RunningTime = 1212.413; Preprocess = $preprocessA; ModelDetails = $modelDetailsA; Test = 1;                 Val = 0.900990099009901; J = 0.0714285714285714 },
    @{ Row = 16; Time = '20160406_120255'; RunningTime = 1304.341; Preprocess = $preprocessA; ModelDetails = $modelDetailsA; Test = 1;                 Val = 0.897689768976898; J = 0.0869565217391304 },
    @{ Row = 17; Time = '20160406_133631'; RunningTime = 2582.359; Preprocess = $preprocessB; ModelDetails = $modelDetailsB; Test = 0.999333333333333; Val = 0.887788778877888; J = 0.0909090909090909 },
    @{ Row = 18; Time = '20160406_141933'; RunningTime = 2581.266; Preprocess = $preprocessB; ModelDetails = $modelDetailsB; Test = 1;                 Val = 0.897689768976898; J = 0.115942028985507 },
    @{ Row = 19; Time = '20160406_150234'; RunningTime = 1689.033; Preprocess = $preprocessB; ModelDetails = $modelDetailsB; Test = 1;                 Val = 0.900990099009901; J = 0.0571428571428571 },
    @{ Row = 20; Time = '20160406_153043'; RunningTime = 1659.883; Preprocess = $preprocessB; ModelDetails = $modelDetailsB; Test = 0.999333333333333; Val = 0.894389438943894; J = 0.102941176470588 },
    @{ Row = 21; Time = '20160406_155823'; RunningTime = 1483.621; Preprocess = $preprocessB; ModelDetails = $modelDetailsB; Test = 0.999333333333333; Val = 0.881188118811881; J = 0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Time
    $ws.Cells.Item($row, 2).Value = $r.RunningTime
    $ws.Cells.Item($row, 3).Value = $r.Preprocess
    $ws.Cells.Item($row, 4).Value = $features
    $ws.Cells.Item($row, 5).Value = $model
    $ws.Cells.Item($row, 6).Value = $r.ModelDetails
    $ws.Cells.Item($row, 7).Value = $r.Test
    $ws.Cells.Item($row, 8).Value = $r.Val
    $ws.Cells.Item($row, 9).Value = $templateFilter
    $ws.Cells.Item($row, 10).Value = $r.J
}
